$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the "Femacal de La Calera - Plátano" data block
# (new weekly price entries), pushing the existing historical rows down by 3.
$ws.Rows("502:504").Insert()

$data = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44509, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Maduro", 160, 17000, 17000, 17000, "$/caja 20 kilos", "Ecuador", 850, 20),
    @(3, "Femacal de La Calera", "Coquimbo", 44509, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Pintón", 240, 18000, 18000, 18000, "$/caja 20 kilos", "Ecuador", 900, 20),
    @(3, "Femacal de La Calera", "Coquimbo", 44509, 5, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Primera Pintón", 240, 20000, 20000, 20000, "$/caja 20 kilos", "Ecuador", 1000, 20)
)

for ($i = 0; $i -lt 3; $i++) {
    $r = 502 + $i
    $row = $data[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
